$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88; this shifts current rows 88-109 down to 89-110,
# carrying their values/formatting with them.
$ws.Rows("88:88").Insert()

# Populate the newly inserted row 88 with the new weekly record.
$ws.Cells.Item(88, 1).Value = 8
$ws.Cells.Item(88, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(88, 3).Value = "Coquimbo"
$ws.Cells.Item(88, 4).Value = 44785
$ws.Cells.Item(88, 5).Value = 4
$ws.Cells.Item(88, 6).Value = 100112052
$ws.Cells.Item(88, 7).Value = "Albahaca"
$ws.Cells.Item(88, 8).Value = "Sin especificar"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 1200
$ws.Cells.Item(88, 11).Value = 3300
$ws.Cells.Item(88, 12).Value = 3500
$ws.Cells.Item(88, 13).Value = 3400
$ws.Cells.Item(88, 14).Value = "$/paquete"
$ws.Cells.Item(88, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(88, 16).Value = 3400
$ws.Cells.Item(88, 17).Value = 1
$ws.Cells.Item(88, 18).Value = "Hortaliza"
